# Set one of the applications to contested:
# L4 ("Application Contested" column for row 4) changes from "No" to "Yes".
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("L4").Value = "Yes"

# Reflect the reviewer's resulting view/selection state: scrolled right so
# column G is at the left edge, with L5 as the active cell.
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 7
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("L5").Select()
